$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 40
$ws_ALC.Range("H40").Value = 4084.4614
$ws_ALC.Range("I40").Value = 2859.9
$ws_ALC.Range("K40").Value = 2859.9
$ws_ALC.Range("M40").Value = -2684.9

# ALC row 41
$ws_ALC.Range("H41").Value = 472.36365
$ws_ALC.Range("I41").Value = 321.77777
$ws_ALC.Range("K41").Value = 321.77777
$ws_ALC.Range("M41").Value = 118.22223

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws_ARM.Range("H2").Value = 2192.1667
$ws_ARM.Range("I2").Value = 1106.5555
$ws_ARM.Range("K2").Value = 1106.5555
$ws_ARM.Range("M2").Value = -993.5554999999999

# ARM row 32
$ws_ARM.Range("H32").Value = 3528.524
$ws_ARM.Range("I32").Value = 2597.1177
$ws_ARM.Range("K32").Value = 2597.1177
$ws_ARM.Range("M32").Value = -2310.1177

# ARM row 63
$ws_ARM.Range("H63").Value = 4744.1816
$ws_ARM.Range("I63").Value = 3885.1428
$ws_ARM.Range("J63").Value = 6247.5
$ws_ARM.Range("K63").Value = 3885.1428
$ws_ARM.Range("L63").Value = 6247.5
$ws_ARM.Range("M63").Value = -3199.1428
$ws_ARM.Range("N63").Value = -7619.5

# ARM row 66
$ws_ARM.Range("H66").Value = 4744.1816
$ws_ARM.Range("I66").Value = 3885.1428
$ws_ARM.Range("J66").Value = 6247.5
$ws_ARM.Range("K66").Value = 19425.714
$ws_ARM.Range("L66").Value = 31237.5
$ws_ARM.Range("M66").Value = -15993.714
$ws_ARM.Range("N66").Value = -38101.5

# ARM row 74
$ws_ARM.Range("H74").Value = 3496.75
$ws_ARM.Range("I74").Value = 3496.75
$ws_ARM.Range("K74").Value = 3496.75
$ws_ARM.Range("M74").Value = -2622.75

# ARM row 77
$ws_ARM.Range("H77").Value = 3496.75
$ws_ARM.Range("I77").Value = 3496.75
$ws_ARM.Range("K77").Value = 17483.75
$ws_ARM.Range("M77").Value = -13115.75

# ARM row 88
$ws_ARM.Range("H88").Value = 1836.375
$ws_ARM.Range("I88").Value = 1299.1111
$ws_ARM.Range("K88").Value = 1299.1111
$ws_ARM.Range("M88").Value = -893.1111000000001

# ARM row 91
$ws_ARM.Range("H91").Value = 1836.375
$ws_ARM.Range("I91").Value = 1299.1111
$ws_ARM.Range("K91").Value = 1299.1111
$ws_ARM.Range("M91").Value = 104.8888999999999

# ARM row 97
$ws_ARM.Range("H97").Value = 1351.9584
$ws_ARM.Range("I97").Value = 1007.3158
$ws_ARM.Range("J97").Value = 2661.6
$ws_ARM.Range("K97").Value = 1007.3158
$ws_ARM.Range("L97").Value = 2661.6
$ws_ARM.Range("M97").Value = -511.3158
$ws_ARM.Range("N97").Value = -3653.6

# ARM row 116
$ws_ARM.Range("H116").Value = 2192.1667
$ws_ARM.Range("I116").Value = 1106.5555
$ws_ARM.Range("K116").Value = 1106.5555
$ws_ARM.Range("M116").Value = 1187.4445

# ARM row 132
$ws_ARM.Range("H132").Value = 1807.25
$ws_ARM.Range("I132").Value = 1661.0667
$ws_ARM.Range("K132").Value = 4983.2001
$ws_ARM.Range("M132").Value = -2453.2001

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws_BSM.Range("H3").Value = 2192.1667
$ws_BSM.Range("I3").Value = 1106.5555
$ws_BSM.Range("K3").Value = 1106.5555
$ws_BSM.Range("M3").Value = -992.5554999999999

# BSM row 82
$ws_BSM.Range("H82").Value = 1955
$ws_BSM.Range("I82").Value = 1955
$ws_BSM.Range("K82").Value = 1955
$ws_BSM.Range("M82").Value = -1572

# BSM row 85
$ws_BSM.Range("H85").Value = 1955
$ws_BSM.Range("I85").Value = 1955
$ws_BSM.Range("K85").Value = 1955
$ws_BSM.Range("M85").Value = -629

# BSM row 86
$ws_BSM.Range("H86").Value = 1665
$ws_BSM.Range("I86").Value = 1926.5714
$ws_BSM.Range("K86").Value = 1926.5714
$ws_BSM.Range("M86").Value = -803.5714

# BSM row 89
$ws_BSM.Range("H89").Value = 1665
$ws_BSM.Range("I89").Value = 1926.5714
$ws_BSM.Range("K89").Value = 9632.857
$ws_BSM.Range("M89").Value = -4016.857

# BSM row 94
$ws_BSM.Range("H94").Value = 1886.3636
$ws_BSM.Range("I94").Value = 2138.889
$ws_BSM.Range("K94").Value = 2138.889
$ws_BSM.Range("M94").Value = -1687.889

# BSM row 97
$ws_BSM.Range("H97").Value = 11100
$ws_BSM.Range("I97").Value = 11100
$ws_BSM.Range("K97").Value = 11100
$ws_BSM.Range("M97").Value = -10109

# BSM row 99
$ws_BSM.Range("H99").Value = 2434.4375
$ws_BSM.Range("I99").Value = 2211.5715
$ws_BSM.Range("K99").Value = 2211.5715
$ws_BSM.Range("M99").Value = -713.5715

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 26
$ws_CRP.Range("H26").Value = 0
$ws_CRP.Range("I26").Value = 0
$ws_CRP.Range("K26").Value = 0
$ws_CRP.Range("M26").ClearContents()

# CRP row 70
$ws_CRP.Range("H70").Value = 25000
$ws_CRP.Range("J70").Value = 25000
$ws_CRP.Range("L70").Value = 25000
$ws_CRP.Range("N70").Value = -25630

# CRP row 73
$ws_CRP.Range("H73").Value = 25000
$ws_CRP.Range("J73").Value = 25000
$ws_CRP.Range("L73").Value = 25000
$ws_CRP.Range("N73").Value = -27184

# CRP row 134
$ws_CRP.Range("H134").Value = 1690.1818
$ws_CRP.Range("I134").Value = 1710.7778
$ws_CRP.Range("J134").Value = 1597.5
$ws_CRP.Range("K134").Value = 5132.3334
$ws_CRP.Range("L134").Value = 4792.5
$ws_CRP.Range("M134").Value = -2597.3334
$ws_CRP.Range("N134").Value = -9862.5

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 32
$ws_CUL.Range("H32").Value = 5000
$ws_CUL.Range("I32").Value = 0
$ws_CUL.Range("K32").Value = 0
$ws_CUL.Range("M32").ClearContents()

# CUL row 64
$ws_CUL.Range("H64").Value = 1100
$ws_CUL.Range("I64").Value = 1100
$ws_CUL.Range("K64").Value = 3300
$ws_CUL.Range("M64").Value = -3030

# CUL row 67
$ws_CUL.Range("H67").Value = 1100
$ws_CUL.Range("I67").Value = 1100
$ws_CUL.Range("K67").Value = 3300
$ws_CUL.Range("M67").Value = -2364

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 46
$ws_GSM.Range("H46").Value = 8165.4443
$ws_GSM.Range("J46").Value = 11665
$ws_GSM.Range("L46").Value = 11665
$ws_GSM.Range("N46").Value = -11977

# GSM row 80
$ws_GSM.Range("H80").Value = 4558
$ws_GSM.Range("J80").Value = 11003
$ws_GSM.Range("L80").Value = 11003
$ws_GSM.Range("N80").Value = -12999

# GSM row 83
$ws_GSM.Range("H83").Value = 4558
$ws_GSM.Range("J83").Value = 11003
$ws_GSM.Range("L83").Value = 55015
$ws_GSM.Range("N83").Value = -64999

# GSM row 97
$ws_GSM.Range("H97").Value = 1957.25
$ws_GSM.Range("I97").Value = 2090
$ws_GSM.Range("J97").Value = 1824.5
$ws_GSM.Range("K97").Value = 2090
$ws_GSM.Range("L97").Value = 1824.5
$ws_GSM.Range("M97").Value = -1594
$ws_GSM.Range("N97").Value = -2816.5

# GSM row 113
$ws_GSM.Range("H113").Value = 2661.3572
$ws_GSM.Range("I113").Value = 2132.1
$ws_GSM.Range("K113").Value = 2132.1
$ws_GSM.Range("M113").Value = 37.90000000000009

# GSM row 132
$ws_GSM.Range("H132").Value = 2087.375
$ws_GSM.Range("I132").Value = 2106.5
$ws_GSM.Range("J132").Value = 2049.125
$ws_GSM.Range("K132").Value = 6319.5
$ws_GSM.Range("L132").Value = 6147.375
$ws_GSM.Range("M132").Value = -3789.5
$ws_GSM.Range("N132").Value = -11207.375

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws_LTW.Range("H7").Value = 4372.7144
$ws_LTW.Range("I7").Value = 4721.8
$ws_LTW.Range("J7").Value = 3500
$ws_LTW.Range("K7").Value = 4721.8
$ws_LTW.Range("L7").Value = 3500
$ws_LTW.Range("M7").Value = -4609.8
$ws_LTW.Range("N7").Value = -3724

# LTW row 22
$ws_LTW.Range("H22").Value = 4199.6
$ws_LTW.Range("J22").Value = 4199.6
$ws_LTW.Range("L22").Value = 4199.6
$ws_LTW.Range("N22").Value = -4789.6

# LTW row 27
$ws_LTW.Range("H27").Value = 4199.6
$ws_LTW.Range("J27").Value = 4199.6
$ws_LTW.Range("L27").Value = 4199.6
$ws_LTW.Range("N27").Value = -4413.6

# LTW row 55
$ws_LTW.Range("H55").Value = 340.23077
$ws_LTW.Range("I55").Value = 120
$ws_LTW.Range("J55").Value = 438.1111
$ws_LTW.Range("K55").Value = 120
$ws_LTW.Range("L55").Value = 438.1111
$ws_LTW.Range("M55").Value = 53
$ws_LTW.Range("N55").Value = -784.1111000000001

# LTW row 82
$ws_LTW.Range("H82").Value = 529.3333
$ws_LTW.Range("I82").Value = 529.3333
$ws_LTW.Range("J82").Value = 0
$ws_LTW.Range("K82").Value = 529.3333
$ws_LTW.Range("L82").Value = 0
$ws_LTW.Range("M82").Value = -168.3333
$ws_LTW.Range("N82").ClearContents()

# LTW row 85
$ws_LTW.Range("H85").Value = 529.3333
$ws_LTW.Range("I85").Value = 529.3333
$ws_LTW.Range("J85").Value = 0
$ws_LTW.Range("K85").Value = 529.3333
$ws_LTW.Range("L85").Value = 0
$ws_LTW.Range("M85").Value = 718.6667
$ws_LTW.Range("N85").ClearContents()

# LTW row 93
$ws_LTW.Range("H93").Value = 753.8
$ws_LTW.Range("I93").Value = 692.25
$ws_LTW.Range("J93").Value = 1000
$ws_LTW.Range("K93").Value = 692.25
$ws_LTW.Range("L93").Value = 1000
$ws_LTW.Range("M93").Value = 555.75
$ws_LTW.Range("N93").Value = -3496

# LTW row 100
$ws_LTW.Range("H100").Value = 15109.6
$ws_LTW.Range("I100").Value = 10337
$ws_LTW.Range("K100").Value = 10337
$ws_LTW.Range("M100").Value = -9796

# LTW row 126
$ws_LTW.Range("H126").Value = 4372.7144
$ws_LTW.Range("I126").Value = 4721.8
$ws_LTW.Range("J126").Value = 3500
$ws_LTW.Range("K126").Value = 14165.4
$ws_LTW.Range("L126").Value = 10500
$ws_LTW.Range("M126").Value = -11695.4
$ws_LTW.Range("N126").Value = -15440

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 93
$ws_WVR.Range("H93").Value = 0
$ws_WVR.Range("J93").Value = 0
$ws_WVR.Range("L93").Value = 0
$ws_WVR.Range("N93").ClearContents()

# WVR row 96
$ws_WVR.Range("H96").Value = 1616.8823
$ws_WVR.Range("I96").Value = 1644
$ws_WVR.Range("J96").Value = 1551.8
$ws_WVR.Range("K96").Value = 1644
$ws_WVR.Range("L96").Value = 1551.8
$ws_WVR.Range("M96").Value = -271
$ws_WVR.Range("N96").Value = -4297.8

# WVR row 113
$ws_WVR.Range("H113").Value = 770.3125
$ws_WVR.Range("I113").Value = 1291.6666
$ws_WVR.Range("J113").Value = 457.5
$ws_WVR.Range("K113").Value = 3874.9998
$ws_WVR.Range("L113").Value = 1372.5
$ws_WVR.Range("M113").Value = -1704.9998
$ws_WVR.Range("N113").Value = -5712.5

# WVR row 136
$ws_WVR.Range("H136").Value = 924.1429000000001
$ws_WVR.Range("I136").Value = 956.7692
$ws_WVR.Range("J136").Value = 500
$ws_WVR.Range("K136").Value = 2870.3076
$ws_WVR.Range("L136").Value = 1500
$ws_WVR.Range("M136").Value = -320.3076000000001
$ws_WVR.Range("N136").Value = -6600
